# Updated cryptos list values (Price / Volume(1h) columns)
# D column values are force-written as text (leading-apostrophe Formula trick)
# so numeric-looking strings (e.g. "229.57") do not get silently converted
# to Excel numbers -- matching the workbook's original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'38.683.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.45%  "

$ws.Range("D3").Formula = "'2.103.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.55%  "

$ws.Range("D5").Formula = "'229.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "

$ws.Range("D6").Formula = "'0.617"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.39%  "

$ws.Range("D7").Formula = "'61.63"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.95%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +2.03%  "

$ws.Range("D10").Formula = "'0.0845"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.31%  "

$ws.Range("D11").Formula = "'0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.89%  "

$ws.Range("D12").Formula = "'2.413.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.49%  "

$ws.Range("D13").Formula = "'14.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.24%  "

$ws.Range("D14").Formula = "'22.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.90%  "

$ws.Range("E15").Value = "  +2.40%  "

$ws.Range("E16").Value = "  +5.82%  "

$ws.Range("D17").Formula = "'2.095.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.47%  "

$ws.Range("D18").Formula = "'38.560.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.25%  "

$ws.Range("D19").Formula = "'71.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.33%  "

$ws.Range("D20").Formula = "'6.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.95%  "

$ws.Range("D21").Formula = "'0.0₃0837"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.41%  "

$ws.Range("D22").Formula = "'227.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.46%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Formula = "'2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.10%  "

$ws.Range("D25").Formula = "'2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.63%  "

$ws.Range("D26").Formula = "'170.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.30%  "

$ws.Range("D27").Formula = "'9.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.04%  "

$ws.Range("E28").Value = "  +1.48%  "

$ws.Range("D29").Formula = "'19.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.06%  "

$ws.Range("D30").Formula = "'1.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.50%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").Formula = "'2.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.50%  "

$ws.Range("D33").Formula = "'4.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.15%  "

$ws.Range("E34").Value = "  +2.44%  "

$ws.Range("D35").Formula = "'0.0608"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("D36").Formula = "'6.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.75%  "

$ws.Range("D37").Formula = "'2.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.99%  "

$ws.Range("D38").Formula = "'3.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.21%  "

$ws.Range("E39").Value = "  -0.16%  "

$ws.Range("D40").Formula = "'18.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.58%  "

$ws.Range("D41").Formula = "'1.547.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.69%  "

$ws.Range("D42").Formula = "'100.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.74%  "

$ws.Range("D43").Formula = "'0.0222"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.91%  "

$ws.Range("E44").Value = "  +1.26%  "

$ws.Range("E45").Value = "  +1.19%  "

$ws.Range("D46").Formula = "'4.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.71%  "

$ws.Range("D47").Formula = "'7.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.42%  "

$ws.Range("E48").Value = "  +1.53%  "

$ws.Range("E49").Value = "  +3.87%  "

$ws.Range("E50").Value = "  +0.59%  "

$ws.Range("D51").Formula = "'2.299.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.52%  "
